# Update the "Page" sheet's destinations / nav-graph data to better match
# the structure of the individual modules.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page")

# The "Start out right..." intro page (row 7) moved up into the module's
# own sheet/flow and is no longer a standalone Page row here - remove it,
# shifting every row below it up by one.
$ws.Rows(7).Delete()

# The module's short description text on row 3 now reuses the "The
# Foundations" string instead of the old placeholder text.
$ws.Range("E3").Value = "The Foundations"

# Renumber the pages that used to follow the removed intro page (now
# rows 7-13) so the page index column (C) is zero-based again.
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 5
$ws.Range("C13").Value = 6

# Restore the cursor/selection to where it was left after the edit.
$ws.Range("C9").Select()
